# Create Interaction System and add interaction UI.
# Adds new localized strings for the L1/Dialogue-1 interaction lines plus a
# couple of UI/common strings (Thanks/Welcome/Yes/No, Cook), wiring them into
# the "Narrative " sheet (new rows 8-25) and the "UI" sheet (new row 68).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "UI"
$ws2 = $wb.Worksheets.Item(2)   # "Narrative "

$ws2.Range("A8").Value = 'STR_L1_Thanks'
$ws2.Range("B8").Value = 'Thanks'
$ws2.Range("C8").Value = 'Cảm ơn.'
$ws2.Range("A9").Value = 'STR_L1_Welcom'
$ws2.Range("B9").Value = 'You are welcome!'
$ws2.Range("A10").Value = 'STR_L1_No'
$ws2.Range("B10").Value = 'No'
$ws2.Range("C10").Value = 'Không.'
$ws2.Range("C9").Value = 'Không có chi.'
$ws2.Range("A11").Value = 'STR_L1_Yes'
$ws2.Range("B11").Value = 'Yes'
$ws2.Range("C11").Value = 'Vâng'
$ws2.Range("A12").Value = 'L1_D1_Default_Jade_Emperor'
$ws2.Range("C12").Value = 'Con khỉ nhà người lại tìm ta có chuyện gì đây?'
$ws2.Range("A13").Value = 'L1_D1_Default_Taibai_Jinxing'
$ws2.Range("C13").Value = 'Không biết Đại Thánh tìm tôi có chuyện gì không?'
$ws2.Range("A14").Value = 'L1_D1_Default_Bodhi_Patriarch'
$ws2.Range("C14").Value = 'Ngộ Không, người muốn học loại phép thuật gì?'
$ws2.Range("A15").Value = 'L1_D1_Default_Taishang_Laojun'
$ws2.Range("C15").Value = 'Không biết Đại Thánh tới đây có việc gì không?'
$ws2.Range("A16").Value = 'L1_D1_Default_Third_Prince_Nezha'
$ws2.Range("C16").Value = 'Xin chào Tôn Đại Thánh.'
$ws2.Range("A17").Value = 'L1_D1_Default_Bull_Demon_ King'
$ws2.Range("C17").Value = 'Người huynh đệ dạo này khỏe không?'
$ws2.Range("A18").Value = 'L1_D1_Default_Dragon_king_Eatern_Sea'
$ws2.Range("C18").Value = 'Ta có thể giúp gì cho ngài?'
$ws2.Range("A19").Value = 'L1_D1_Default_Sha_Wujing'
$ws2.Range("C19").Value = 'Đại sư huynh cứ để đệ trông coi hành lý cho.'
$ws2.Range("A20").Value = 'L1_D1_Default_Zhu_Baije'
$ws2.Range("C20").Value = 'Có phải yêu quá đến phải không đại sư huynh.'
$ws2.Range("A21").Value = 'L1_D1_Default_Tang_SanZanng'
$ws2.Range("C21").Value = 'Ngộ không, không được vô lễ.'
$ws2.Range("A22").Value = 'L1_D1_Default_Little_White_Dragon'
$ws2.Range("C22").Value = '"Hí hí hí"'
$ws2.Range("A23").Value = 'L1_D1_Default_Guanyin_Bodhisattva'
$ws2.Range("C23").Value = 'Ngộ không sao ngươi không đi bảo vệ sự phụ lại đến đây tìm ta có chuyện gì không?'
$ws2.Range("A24").Value = 'L1_D1_Default_Shakyamun_ Buddha'
$ws2.Range("C24").Value = 'Con khí này lại đến rồi.'
$ws2.Range("A25").Value = 'L1_D1_Default_Yanluo_Wang'
$ws2.Range("C25").Value = 'Đại Thánh tìm Tiểu Vương có việc gì thế?'
$ws1.Range("A68").Value = 'UI_COOK'
$ws1.Range("B68").Value = 'Cook'
$ws1.Range("C68").Value = 'Nấu Nướng'

# Widen the ID column on the Narrative sheet now that longer interaction IDs
# (e.g. L1_D1_Default_Dragon_king_Eatern_Sea) live there.
$ws2.Columns.Item(1).ColumnWidth = 31.7

# Selection / active-sheet bookkeeping: make the Narrative sheet's selection
# land on the newly added block, then switch back to the UI sheet (now the
# active tab) with its selection on the freshly appended row.
$ws2.Range("C8").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("E65").Select() | Out-Null
